# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets.
# These two sheets carry duplicated event data; both need the same updates.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 235
    3  = 263
    5  = 815
    6  = 269
    7  = 6511
    12 = 35
    14 = 9
    15 = 207
    16 = 523
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
